$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first 9 shops (rows 2-10) get their ItemID (column C) updated from the
# generic "BuildN" placeholders to the new named SLG buildings.
$ws.Range("C2").Value = "Build_Altar_1"
$ws.Range("C3").Value = "Build_Arena_1"
$ws.Range("C4").Value = "Build_Camp_1"
$ws.Range("C5").Value = "Build_Gold_Mine_1"
$ws.Range("C6").Value = "Build_Item_Hourse_1"
$ws.Range("C7").Value = "Build_League_1"
$ws.Range("C8").Value = "Build_Magic_Hourse_1"
$ws.Range("C9").Value = "Build_Tower_1"
$ws.Range("C10").Value = "Build_Town_1"

[void]$ws.Range("D11").Select()
